$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2: Jason Poage -> Rachel Tipton
$ws.Range("A2").Value = "Rachel"
$ws.Range("B2").Value = "Tipton"
$ws.Range("D2").Value = "316 Louisville - Retail"
$ws.Range("E2").Value = "Operation Manager"
$ws.Range("F2").Value = "Shannon Drabant"
$ws.Range("G2").Value = 13901022
$ws.Range("H2").Value = "rntipton@charter.net"
$ws.Range("J2").Value = "2021-12-30T00:00:00"
$ws.Range("K2").Value = "2021-12-30T00:00:00"
$ws.Range("N2").Value = "2022-01-01T17:04:20.363"

# New row 3: Dirk Tomlinson
$ws.Range("A3").Value = "Dirk"
$ws.Range("B3").Value = "Tomlinson"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = "301 Evansville - Retail"
$ws.Range("E3").Value = "Technician"
$ws.Range("F3").Value = "Dakota Floyd"
$ws.Range("G3").Value = 13901023
$ws.Range("H3").Value = "dtomlinson145@gmail.com"
$ws.Range("I3").Value = $true
$ws.Range("J3").Value = "2022-01-04T00:00:00"
$ws.Range("K3").Value = "2022-01-02T00:00:00"
$ws.Range("L3").Value = "0001-01-01T00:00:00"
$ws.Range("M3").Value = "POB Completed"
$ws.Range("N3").Value = "2022-01-03T14:39:15.637"
$ws.Range("O3").Value = "0001-01-01T00:00:00"
